$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(4, 6).Value = 1760  # was 1755
$ws.Cells.Item(9, 6).Value = 313  # was 307
$ws.Cells.Item(10, 6).Value = 1700  # was 1698
$ws.Cells.Item(12, 6).Value = 1402  # was 1399
$ws.Cells.Item(14, 6).Value = 322  # was 321
$ws.Cells.Item(15, 6).Value = 664  # was 663
$ws.Cells.Item(16, 6).Value = 12662  # was 12650
$ws.Cells.Item(17, 6).Value = 12685  # was 12672
$ws.Cells.Item(18, 6).Value = 940  # was 937
$ws.Cells.Item(19, 6).Value = 733  # was 732
$ws.Cells.Item(21, 6).Value = 300  # was 299
$ws.Cells.Item(22, 6).Value = 46  # was 45
$ws.Cells.Item(23, 6).Value = 516  # was 509
$ws.Cells.Item(24, 6).Value = 1985  # was 1981
$ws.Cells.Item(27, 6).Value = 232  # was 231
$ws.Cells.Item(28, 6).Value = 664  # was 663

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(5, 6).Value = 75  # was 73
$ws.Cells.Item(6, 6).Value = 12  # was 11
$ws.Cells.Item(9, 6).Value = 51  # was 50
$ws.Cells.Item(11, 6).Value = 3  # was 2

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 154  # was 153

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 154  # was 153
$ws.Cells.Item(6, 6).Value = 1760  # was 1755
$ws.Cells.Item(14, 6).Value = 313  # was 307
$ws.Cells.Item(15, 6).Value = 1700  # was 1698
$ws.Cells.Item(17, 6).Value = 1402  # was 1399
$ws.Cells.Item(19, 6).Value = 322  # was 321
$ws.Cells.Item(20, 6).Value = 75  # was 73
$ws.Cells.Item(21, 6).Value = 664  # was 663
$ws.Cells.Item(22, 6).Value = 12662  # was 12650
$ws.Cells.Item(23, 6).Value = 12685  # was 12672
$ws.Cells.Item(24, 6).Value = 940  # was 937
$ws.Cells.Item(25, 6).Value = 733  # was 732
$ws.Cells.Item(27, 6).Value = 300  # was 299
$ws.Cells.Item(28, 6).Value = 46  # was 45
$ws.Cells.Item(29, 6).Value = 516  # was 509
$ws.Cells.Item(30, 6).Value = 12  # was 11
$ws.Cells.Item(32, 6).Value = 1985  # was 1981
$ws.Cells.Item(36, 6).Value = 51  # was 50
$ws.Cells.Item(37, 6).Value = 232  # was 231
$ws.Cells.Item(38, 6).Value = 664  # was 663
$ws.Cells.Item(40, 6).Value = 3  # was 2
